$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_5_1_23"
$ws.Range("B2").Value = 0.4761419077363901
$ws.Range("C2").Value = -4.095801797538899
$ws.Range("D2").Value = 0.5396644167337645
$ws.Range("E2").Value = -1.263026450368378
$ws.Range("F2").Value = 0.5797566175460815
$ws.Range("G2").Value = 3.013367652893066
$ws.Range("H2").Value = 0.4783340096473694
$ws.Range("I2").Value = 1.82040810585022

$ws.Range("A3").Value = "model_5_1_21"
$ws.Range("B3").Value = 0.4765386897312088
$ws.Range("C3").Value = -4.10249486835295
$ws.Range("D3").Value = 0.5626151518461262
$ws.Range("E3").Value = -1.251681194293142
$ws.Range("F3").Value = 0.5793174505233765
$ws.Range("G3").Value = 3.017325401306152
$ws.Range("H3").Value = 0.4544859230518341
$ws.Range("I3").Value = 1.81128191947937

$ws.Range("A4").Value = "model_5_1_24"
$ws.Range("B4").Value = 0.4766864133483077
$ws.Range("C4").Value = -4.083338424558931
$ws.Range("D4").Value = 0.5285503298378746
$ws.Range("E4").Value = -1.264931455587787
$ws.Range("F4").Value = 0.5791538953781128
$ws.Range("G4").Value = 3.0059974193573
$ws.Range("H4").Value = 0.4898826479911804
$ws.Range("I4").Value = 1.821940422058105

$ws.Range("A5").Value = "model_5_1_20"
$ws.Range("B5").Value = 0.4774240316017642
$ws.Range("C5").Value = -4.101403871988448
$ws.Range("D5").Value = 0.5767565686652105
$ws.Range("E5").Value = -1.242667207831032
$ws.Range("F5").Value = 0.5783376097679138
$ws.Range("G5").Value = 3.016680240631104
$ws.Range("H5").Value = 0.4397916197776794
$ws.Range("I5").Value = 1.804031014442444

$ws.Range("A6").Value = "model_5_1_22"
$ws.Range("B6").Value = 0.4777693409269723
$ws.Range("C6").Value = -4.090401523272354
$ws.Range("D6").Value = 0.5559837331558601
$ws.Range("E6").Value = -1.251003887440886
$ws.Range("F6").Value = 0.577955424785614
$ws.Range("G6").Value = 3.010174036026001
$ws.Range("H6").Value = 0.4613766372203827
$ws.Range("I6").Value = 1.810737013816833

$ws.Range("A7").Value = "model_5_1_19"
$ws.Range("B7").Value = 0.481233901901883
$ws.Range("C7").Value = -4.075485473202187
$ws.Range("D7").Value = 0.5971623707855923
$ws.Range("E7").Value = -1.220171355828374
$ws.Range("F7").Value = 0.5741211771965027
$ws.Range("G7").Value = 3.00135350227356
$ws.Range("H7").Value = 0.4185879826545715
$ws.Range("I7").Value = 1.785935163497925

$ws.Range("A8").Value = "model_5_1_17"
$ws.Range("B8").Value = 0.4824468393807838
$ws.Range("C8").Value = -4.066450278644949
$ws.Range("D8").Value = 0.6184595120635119
$ws.Range("E8").Value = -1.203710734548255
$ws.Range("F8").Value = 0.5727788209915161
$ws.Range("G8").Value = 2.996010780334473
$ws.Range("H8").Value = 0.3964581489562988
$ws.Range("I8").Value = 1.772693753242493

$ws.Range("A9").Value = "model_5_1_18"
$ws.Range("B9").Value = 0.4849467135075098
$ws.Range("C9").Value = -4.05205699845357
$ws.Range("D9").Value = 0.6187793465412391
$ws.Range("E9").Value = -1.197910783048083
$ws.Range("F9").Value = 0.5700122117996216
$ws.Range("G9").Value = 2.987499237060547
$ws.Range("H9").Value = 0.3961258232593536
$ws.Range("I9").Value = 1.768028259277344

$ws.Range("A10").Value = "model_5_1_16"
$ws.Range("B10").Value = 0.4869863644246838
$ws.Range("C10").Value = -4.031796710299143
$ws.Range("D10").Value = 0.6398672675244793
$ws.Range("E10").Value = -1.177208955692927
$ws.Range("F10").Value = 0.567754864692688
$ws.Range("G10").Value = 2.975518226623535
$ws.Range("H10").Value = 0.3742133975028992
$ws.Range("I10").Value = 1.751375436782837

$ws.Range("A11").Value = "model_5_1_14"
$ws.Range("B11").Value = 0.4893515117463689
$ws.Range("C11").Value = -4.00193312499606
$ws.Range("D11").Value = 0.6588720022438006
$ws.Range("E11").Value = -1.154031129753027
$ws.Range("F11").Value = 0.5651373863220215
$ws.Range("G11").Value = 2.957858800888062
$ws.Range("H11").Value = 0.3544656038284302
$ws.Range("I11").Value = 1.732730865478516

$ws.Range("A12").Value = "model_5_1_15"
$ws.Range("B12").Value = 0.4915605992622745
$ws.Range("C12").Value = -3.993904407357566
$ws.Range("D12").Value = 0.6598769432057522
$ws.Range("E12").Value = -1.150296132732791
$ws.Range("F12").Value = 0.5626925826072693
$ws.Range("G12").Value = 2.953111171722412
$ws.Range("H12").Value = 0.3534213602542877
$ws.Range("I12").Value = 1.729726433753967

$ws.Range("A13").Value = "model_5_1_12"
$ws.Range("B13").Value = 0.4951897045390933
$ws.Range("C13").Value = -3.936097288591609
$ws.Range("D13").Value = 0.6820427447052408
$ws.Range("E13").Value = -1.114326434929588
$ws.Range("F13").Value = 0.5586762428283691
$ws.Range("G13").Value = 2.918927431106567
$ws.Range("H13").Value = 0.3303889036178589
$ws.Range("I13").Value = 1.700791954994202

$ws.Range("A14").Value = "model_5_1_13"
$ws.Range("B14").Value = 0.4958184457278003
$ws.Range("C14").Value = -3.944596768503607
$ws.Range("D14").Value = 0.6810935014025693
$ws.Range("E14").Value = -1.118210469596824
$ws.Range("F14").Value = 0.557980477809906
$ws.Range("G14").Value = 2.923953294754028
$ws.Range("H14").Value = 0.3313752412796021
$ws.Range("I14").Value = 1.703916311264038

$ws.Range("A15").Value = "model_5_1_11"
$ws.Range("B15").Value = 0.4963087244026392
$ws.Range("C15").Value = -3.905912166575102
$ws.Range("D15").Value = 0.6826661139233883
$ws.Range("E15").Value = -1.102195859506671
$ws.Range("F15").Value = 0.5574377775192261
$ws.Range("G15").Value = 2.901077508926392
$ws.Range("H15").Value = 0.3297411799430847
$ws.Range("I15").Value = 1.691033840179443

$ws.Range("A16").Value = "model_5_1_10"
$ws.Range("B16").Value = 0.4974854132247175
$ws.Range("C16").Value = -3.879881009805841
$ws.Range("D16").Value = 0.6867885265791033
$ws.Range("E16").Value = -1.089563774922466
$ws.Range("F16").Value = 0.5561355352401733
$ws.Range("G16").Value = 2.885684251785278
$ws.Range("H16").Value = 0.3254575729370117
$ws.Range("I16").Value = 1.680872678756714

$ws.Range("A17").Value = "model_5_1_9"
$ws.Range("B17").Value = 0.503617293770444
$ws.Range("C17").Value = -3.784354494924904
$ws.Range("D17").Value = 0.6965871710566809
$ws.Range("E17").Value = -1.046429815081998
$ws.Range("F17").Value = 0.5493493676185608
$ws.Range("G17").Value = 2.829195261001587
$ws.Range("H17").Value = 0.3152758181095123
$ws.Range("I17").Value = 1.646174907684326

$ws.Range("A18").Value = "model_5_1_8"
$ws.Range("B18").Value = 0.5086602053857217
$ws.Range("C18").Value = -3.647580259934213
$ws.Range("D18").Value = 0.6885611307976636
$ws.Range("E18").Value = -0.998075093298544
$ws.Range("F18").Value = 0.5437683463096619
$ws.Range("G18").Value = 2.74831485748291
$ws.Range("H18").Value = 0.3236156702041626
$ws.Range("I18").Value = 1.607277870178223

$ws.Range("A19").Value = "model_5_1_7"
$ws.Range("B19").Value = 0.5173289973897821
$ws.Range("C19").Value = -3.544159442224899
$ws.Range("D19").Value = 0.7020495284289701
$ws.Range("E19").Value = -0.9496267564050509
$ws.Range("F19").Value = 0.5341745615005493
$ws.Range("G19").Value = 2.68715763092041
$ws.Range("H19").Value = 0.3095998764038086
$ws.Range("I19").Value = 1.568305373191833

$ws.Range("A20").Value = "model_5_1_6"
$ws.Range("B20").Value = 0.5311261434811714
$ws.Range("C20").Value = -3.397742220819999
$ws.Range("D20").Value = 0.7260937488938827
$ws.Range("E20").Value = -0.8780319028297687
$ws.Range("F20").Value = 0.5189051628112793
$ws.Range("G20").Value = 2.600574731826782
$ws.Range("H20").Value = 0.2846155762672424
$ws.Range("I20").Value = 1.510713458061218

$ws.Range("A21").Value = "model_5_1_5"
$ws.Range("B21").Value = 0.5431274173164427
$ws.Range("C21").Value = -3.252826505352473
$ws.Range("D21").Value = 0.7389773314111288
$ws.Range("E21").Value = -0.8138011557926845
$ws.Range("F21").Value = 0.5056233406066895
$ws.Range("G21").Value = 2.514879941940308
$ws.Range("H21").Value = 0.2712282538414001
$ws.Range("I21").Value = 1.45904552936554

$ws.Range("A22").Value = "model_5_1_4"
$ws.Range("B22").Value = 0.5567939767741518
$ws.Range("C22").Value = -3.096475773853744
$ws.Range("D22").Value = 0.7575777931417222
$ws.Range("E22").Value = -0.7416414454960312
$ws.Range("F22").Value = 0.4904984831809998
$ws.Range("G22").Value = 2.422422885894775
$ws.Range("H22").Value = 0.2519005239009857
$ws.Range("I22").Value = 1.400999188423157

$ws.Range("A23").Value = "model_5_1_3"
$ws.Range("B23").Value = 0.5743209130883771
$ws.Range("C23").Value = -2.896958904612423
$ws.Range("D23").Value = 0.7774690283001513
$ws.Range("E23").Value = -0.6519018806129224
$ws.Range("F23").Value = 0.4711012840270996
$ws.Range("G23").Value = 2.304439783096313
$ws.Range("H23").Value = 0.2312316000461578
$ws.Range("I23").Value = 1.328811645507812

$ws.Range("A24").Value = "model_5_1_2"
$ws.Range("B24").Value = 0.5935203787152864
$ws.Range("C24").Value = -2.709629418759154
$ws.Range("D24").Value = 0.8165541097718615
$ws.Range("E24").Value = -0.555238378376955
$ws.Range("F24").Value = 0.4498531520366669
$ws.Range("G24").Value = 2.193663835525513
$ws.Range("H24").Value = 0.1906183511018753
$ws.Range("I24").Value = 1.251054167747498

$ws.Range("A25").Value = "model_5_1_1"
$ws.Range("B25").Value = 0.609044204839011
$ws.Range("C25").Value = -2.509283556251762
$ws.Range("D25").Value = 0.8226779995473522
$ws.Range("E25").Value = -0.4735440454541375
$ws.Range("F25").Value = 0.4326728284358978
$ws.Range("G25").Value = 2.075190544128418
$ws.Range("H25").Value = 0.1842550337314606
$ws.Range("I25").Value = 1.185338139533997

$ws.Range("A26").Value = "model_5_1_0"
$ws.Range("B26").Value = 0.6311419034690847
$ws.Range("C26").Value = -2.312588799739783
$ws.Range("D26").Value = 0.8776237462809896
$ws.Range("E26").Value = -0.3635933109126874
$ws.Range("F26").Value = 0.4082172214984894
$ws.Range("G26").Value = 1.958876729011536
$ws.Range("H26").Value = 0.127160981297493
$ws.Range("I26").Value = 1.096892356872559

